$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A83").Value = "bland"
$ws.Range("A84").Value = "m7"
$ws.Range("A85").Value = "m7b"
$ws.Range("A86").Value = "za_nr1"

$ws.Range("B83").Value = "State name in Integrated Survey Dataset"
$ws.Range("B84").Value = "Survey Party vote choice variable (from 1973 onwards) "
$ws.Range("B85").Value = "Survey Party vote choice variable (until 1970) "
$ws.Range("B86").Value = "ZA Dataset ID"

[void]$ws.Range("B87").Select()
